$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.203.25'
$ws.Range('E2').Value = '  -0.65%  '

$ws.Range('D3').Value = '1.829.78'
$ws.Range('E3').Value = '  -0.74%  '

$ws.Range('E4').Value = '  +0.16%  '

$ws.Range('D5').Value = "'" + '237.31'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.24%  '

$ws.Range('D6').Value = "'" + '0.6093'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.96%  '

$ws.Range('E7').Value = '  +0.19%  '

$ws.Range('D8').Value = "'" + '0.07089'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -5.21%  '

$ws.Range('D9').Value = "'" + '0.2819'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.93%  '

$ws.Range('D10').Value = "'" + '23.84'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.08%  '

$ws.Range('D11').Value = "'" + '0.07640'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.37%  '

$ws.Range('D12').Value = '1.831.53'
$ws.Range('E12').Value = '  -1.44%  '

$ws.Range('D13').Value = "'" + '4.816'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.49%  '

$ws.Range('D14').Value = "'" + '0.6343'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -6.73%  '

$ws.Range('D15').Value = "'" + '0.000009975'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.77%  '

$ws.Range('D16').Value = '2.070.76'
$ws.Range('E16').Value = '  -1.06%  '

$ws.Range('D17').Value = "'" + '79.34'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.27%  '

$ws.Range('D18').Value = "'" + '5.954'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -5.00%  '

$ws.Range('D19').Value = '29.209.76'
$ws.Range('E19').Value = '  -0.61%  '

$ws.Range('D20').Value = "'" + '228.65'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.79%  '

$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').Value = "'" + '11.81'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.33%  '

$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').Value = "'" + '1.001'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.17%  '

$ws.Range('D23').Value = "'" + '7.045'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.15%  '

$ws.Range('D24').Value = "'" + '1.000'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.02%  '

$ws.Range('D25').Value = "'" + '155.52'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.65%  '

$ws.Range('D26').Value = "'" + '8.103'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.78%  '

$ws.Range('D27').Value = "'" + '0.1304'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.19%  '

$ws.Range('D28').Value = "'" + '16.71'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.49%  '

$ws.Range('D29').Value = "'" + '0.06776'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.10%  '

$ws.Range('D30').Value = "'" + '1.482'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.50%  '

$ws.Range('D31').Value = "'" + '1.458'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.12%  '

$ws.Range('D32').Value = "'" + '3.854'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.43%  '

$ws.Range('D33').Value = "'" + '3.839'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -5.29%  '

$ws.Range('E34').Value = '  -0.93%  '

$ws.Range('D35').Value = "'" + '1.740'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.58%  '

$ws.Range('D36').Value = "'" + '0.6567'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.41%  '

$ws.Range('D37').Value = "'" + '2.558'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.78%  '

$ws.Range('D38').Value = '1.235.83'
$ws.Range('E38').Value = '  -1.18%  '

$ws.Range('D39').Value = "'" + '2.760'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.09%  '

$ws.Range('D40').Value = "'" + '0.01764'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.12%  '

$ws.Range('E41').Value = '  -2.69%  '

$ws.Range('D42').Value = "'" + '0.9237'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.63%  '

$ws.Range('E44').Value = '  -0.74%  '

$ws.Range('D45').Value = "'" + '100.95'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.20%  '

$ws.Range('D46').Value = "'" + '63.68'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.71%  '

$ws.Range('D47').Value = "'" + '0.00000000115'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.94%  '

$ws.Range('D48').Value = "'" + '1.628'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.42%  '

$ws.Range('D49').Value = "'" + '8.588'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -5.18%  '

$ws.Range('D50').Value = "'" + '6.537'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -7.58%  '

$ws.Range('D51').Value = "'" + '0.1084'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.72%  '
